# Add the new "dobno_spolna" worksheet as the last tab (after "Odseljeni"),
# matching the workbook's new 12th sheet.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "dobno_spolna"


# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Županija"
$ws.Range("B1").Value = "Jedinica lokalne samouprave"
$ws.Range("C1").Value = "County of"
$ws.Range("D1").Value = "Local self-government unit"
$ws.Range("E1").Value = "Grad/općina  `nTown/Municipality "
$ws.Range("F1").Value = "Naselje`nSettlement "
$ws.Range("G1").Value = "Spol"
$ws.Range("H1").Value = "Sex"
$ws.Range("I1").Value = "Ukupno`nTotal"
$ws.Range("J1").Value = "0 – 4"
$ws.Range("K1").Value = "5 – 9"
$ws.Range("L1").Value = "10 – 14"
$ws.Range("M1").Value = "15 – 19"
$ws.Range("N1").Value = "20 – 24"
$ws.Range("O1").Value = "25 – 29"
$ws.Range("P1").Value = "30 – 34"
$ws.Range("Q1").Value = "35 – 39"
$ws.Range("R1").Value = "40 – 44"
$ws.Range("S1").Value = "45 – 49"
$ws.Range("T1").Value = "50 – 54"
$ws.Range("U1").Value = "55 – 59"
$ws.Range("V1").Value = "60 – 64"
$ws.Range("W1").Value = "65 – 69"
$ws.Range("X1").Value = "70 – 74"
$ws.Range("Y1").Value = "75 – 79"
$ws.Range("Z1").Value = "80 – 84"
$ws.Range("AA1").Value = "85 – 89"
$ws.Range("AB1").Value = "90 – 94"
$ws.Range("AC1").Value = ">"
$ws.Range("E1").WrapText = $true
$ws.Range("F1").WrapText = $true
$ws.Range("I1").WrapText = $true
$ws.Range("AC1").WrapText = $true
# ---- Row 2 ----
$ws.Range("A2").Value = "Republika Hrvatska"
$ws.Range("C2").Value = "Republic of Croatia"
$ws.Range("G2").Value = "sv."
$ws.Range("H2").Value = "All"
$ws.Range("I2").Value = 3871833
$ws.Range("I2").NumberFormat = "#,##0"
$ws.Range("J2").Value = 175535
$ws.Range("J2").NumberFormat = "#,##0"
$ws.Range("K2").Value = 181445
$ws.Range("K2").NumberFormat = "#,##0"
$ws.Range("L2").Value = 195436
$ws.Range("L2").NumberFormat = "#,##0"
$ws.Range("M2").Value = 188729
$ws.Range("M2").NumberFormat = "#,##0"
$ws.Range("N2").Value = 208852
$ws.Range("N2").NumberFormat = "#,##0"
$ws.Range("O2").Value = 214023
$ws.Range("O2").NumberFormat = "#,##0"
$ws.Range("P2").Value = 227551
$ws.Range("P2").NumberFormat = "#,##0"
$ws.Range("Q2").Value = 255617
$ws.Range("Q2").NumberFormat = "#,##0"
$ws.Range("R2").Value = 267349
$ws.Range("R2").NumberFormat = "#,##0"
$ws.Range("S2").Value = 260146
$ws.Range("S2").NumberFormat = "#,##0"
$ws.Range("T2").Value = 260056
$ws.Range("T2").NumberFormat = "#,##0"
$ws.Range("U2").Value = 279504
$ws.Range("U2").NumberFormat = "#,##0"
$ws.Range("V2").Value = 288351
$ws.Range("V2").NumberFormat = "#,##0"
$ws.Range("W2").Value = 279106
$ws.Range("W2").NumberFormat = "#,##0"
$ws.Range("X2").Value = 228612
$ws.Range("X2").NumberFormat = "#,##0"
$ws.Range("Y2").Value = 146855
$ws.Range("Y2").NumberFormat = "#,##0"
$ws.Range("Z2").Value = 122719
$ws.Range("Z2").NumberFormat = "#,##0"
$ws.Range("AA2").Value = 67249
$ws.Range("AA2").NumberFormat = "#,##0"
$ws.Range("AB2").Value = 21019
$ws.Range("AB2").NumberFormat = "#,##0"
$ws.Range("AC2").Value = 3679
$ws.Range("AC2").NumberFormat = "#,##0"
# ---- Row 3 ----
$ws.Range("A3").Value = "Republika Hrvatska"
$ws.Range("C3").Value = "Republic of Croatia"
$ws.Range("G3").Value = "m"
$ws.Range("H3").Value = "M"
$ws.Range("I3").Value = 1865129
$ws.Range("I3").NumberFormat = "#,##0"
$ws.Range("J3").Value = 90245
$ws.Range("J3").NumberFormat = "#,##0"
$ws.Range("K3").Value = 93311
$ws.Range("K3").NumberFormat = "#,##0"
$ws.Range("L3").Value = 100216
$ws.Range("L3").NumberFormat = "#,##0"
$ws.Range("M3").Value = 97228
$ws.Range("M3").NumberFormat = "#,##0"
$ws.Range("N3").Value = 107102
$ws.Range("N3").NumberFormat = "#,##0"
$ws.Range("O3").Value = 109139
$ws.Range("O3").NumberFormat = "#,##0"
$ws.Range("P3").Value = 114778
$ws.Range("P3").NumberFormat = "#,##0"
$ws.Range("Q3").Value = 128398
$ws.Range("Q3").NumberFormat = "#,##0"
$ws.Range("R3").Value = 134213
$ws.Range("R3").NumberFormat = "#,##0"
$ws.Range("S3").Value = 130035
$ws.Range("S3").NumberFormat = "#,##0"
$ws.Range("T3").Value = 127953
$ws.Range("T3").NumberFormat = "#,##0"
$ws.Range("U3").Value = 134655
$ws.Range("U3").NumberFormat = "#,##0"
$ws.Range("V3").Value = 136338
$ws.Range("V3").NumberFormat = "#,##0"
$ws.Range("W3").Value = 129728
$ws.Range("W3").NumberFormat = "#,##0"
$ws.Range("X3").Value = 100506
$ws.Range("X3").NumberFormat = "#,##0"
$ws.Range("Y3").Value = 59065
$ws.Range("Y3").NumberFormat = "#,##0"
$ws.Range("Z3").Value = 44672
$ws.Range("Z3").NumberFormat = "#,##0"
$ws.Range("AA3").Value = 21206
$ws.Range("AA3").NumberFormat = "#,##0"
$ws.Range("AB3").Value = 5604
$ws.Range("AB3").NumberFormat = "#,##0"
$ws.Range("AC3").Value = 737
# ---- Row 4 ----
$ws.Range("A4").Value = "Republika Hrvatska"
$ws.Range("C4").Value = "Republic of Croatia"
$ws.Range("G4").Value = "ž"
$ws.Range("H4").Value = "W"
$ws.Range("I4").Value = 2006704
$ws.Range("I4").NumberFormat = "#,##0"
$ws.Range("J4").Value = 85290
$ws.Range("J4").NumberFormat = "#,##0"
$ws.Range("K4").Value = 88134
$ws.Range("K4").NumberFormat = "#,##0"
$ws.Range("L4").Value = 95220
$ws.Range("L4").NumberFormat = "#,##0"
$ws.Range("M4").Value = 91501
$ws.Range("M4").NumberFormat = "#,##0"
$ws.Range("N4").Value = 101750
$ws.Range("N4").NumberFormat = "#,##0"
$ws.Range("O4").Value = 104884
$ws.Range("O4").NumberFormat = "#,##0"
$ws.Range("P4").Value = 112773
$ws.Range("P4").NumberFormat = "#,##0"
$ws.Range("Q4").Value = 127219
$ws.Range("Q4").NumberFormat = "#,##0"
$ws.Range("R4").Value = 133136
$ws.Range("R4").NumberFormat = "#,##0"
$ws.Range("S4").Value = 130111
$ws.Range("S4").NumberFormat = "#,##0"
$ws.Range("T4").Value = 132103
$ws.Range("T4").NumberFormat = "#,##0"
$ws.Range("U4").Value = 144849
$ws.Range("U4").NumberFormat = "#,##0"
$ws.Range("V4").Value = 152013
$ws.Range("V4").NumberFormat = "#,##0"
$ws.Range("W4").Value = 149378
$ws.Range("W4").NumberFormat = "#,##0"
$ws.Range("X4").Value = 128106
$ws.Range("X4").NumberFormat = "#,##0"
$ws.Range("Y4").Value = 87790
$ws.Range("Y4").NumberFormat = "#,##0"
$ws.Range("Z4").Value = 78047
$ws.Range("Z4").NumberFormat = "#,##0"
$ws.Range("AA4").Value = 46043
$ws.Range("AA4").NumberFormat = "#,##0"
$ws.Range("AB4").Value = 15415
$ws.Range("AB4").NumberFormat = "#,##0"
$ws.Range("AC4").Value = 2942
$ws.Range("AC4").NumberFormat = "#,##0"


# Match the recorded selection on the new sheet.
$null = $ws.Range("K7").Select()
